$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Version bump: 1.0 -> 1.2.5 (D2 holds the version value; C2 is just the "Version: " label)
$ws.Range("D2").Value = "1.2.5"

# Precondition text: accent + trailing period added.
# The same shared string is used for every test case's precondition (B8, B16, B23, B30, B37).
$newPrecondition = "O usuário devidamente autenticado e na tela inicial do sistema."
$ws.Range("B8").Value = $newPrecondition
$ws.Range("B16").Value = $newPrecondition
$ws.Range("B23").Value = $newPrecondition
$ws.Range("B30").Value = $newPrecondition
$ws.Range("B37").Value = $newPrecondition

# TC1 step text: reworded + accent
$ws.Range("B10").Value = "Beneficiário Acessa o caso de uso através do menu."

# TC2 block: now holds the "detalhar diária" content (was TC3's content)
$ws.Range("B18").Value = "Beneficiário Clica em detalhar diária."
$ws.Range("D18").Value = "SYSTEM Apresenta a tela de Detalhar Diárias."

# TC3 block: now holds the "cancelar diária" content (was TC4's content)
$ws.Range("B25").Value = "Beneficiário Clica em cancelar diária."
$ws.Range("D25").Value = "SYSTEM Apresenta a tela de Cancelar Solicitação de Diária."

# TC4 block: now holds the "analisar prestação de contas" content (was TC2's content)
$ws.Range("B32").Value = "Beneficiário Clica em analisar prestação de contas."
$ws.Range("D32").Value = "SYSTEM Apresenta a tela de Analisar Prestação de Contas."
